$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.650.17"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "3.281.38"
$ws.Range("E3").Value = "  -0.84%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "183.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "3.279.09"
$ws.Range("E8").Value = "  -0.81%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.570"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.45%  "
$ws.Range("E10").Value = "  -5.32%  "
$ws.Range("E11").Value = "  -2.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.19"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.83%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000263"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "3.798.91"
$ws.Range("E14").Value = "  -1.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.39"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "613.89"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.36%  "
$ws.Range("D17").Value = "65.700.51"
$ws.Range("E17").Value = "  -0.55%  "
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.03%  "
$ws.Range("D20").Value = "3.273.39"
$ws.Range("E20").Value = "  -0.88%  "
$ws.Range("E21").Value = "  -3.36%  "
$ws.Range("E22").Value = "  -2.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "98.76"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.00%  "
$ws.Range("E25").Value = "  -0.40%  "
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.43"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.97%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "30.78"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.46"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.70"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "10.82"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "542.42"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.34%  "
$ws.Range("D35").Value = "3.795.43"
$ws.Range("E35").Value = "  -0.99%  "
$ws.Range("E36").Value = "  -2.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "55.96"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.128"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "32.48"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.77%  "
$ws.Range("E41").Value = "  +3.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.13"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.10%  "
$ws.Range("D43").Value = "0.0₃0677"
$ws.Range("E43").Value = "  -7.81%  "
$ws.Range("E44").Value = "  -4.32%  "
$ws.Range("E45").Value = "  -1.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0405"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.17%  "
$ws.Range("E48").Value = "  +0.06%  "
$ws.Range("E49").Value = "  -2.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "127.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.51%  "
